$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Efnb2", "Ephb4", "ECs", 3, 1, 36.899643, 110.698929, 0.7238945645409351, 0.7238945645409351, 2, 0.6666666666666666, 23.746319, 71.238957, 0.7135031414879517, 0.7135031414879517, 876.230693664117, 7886.076242977052, 0.51650104590601, 0.51650104590601)
for ($i=0; $i -lt $row2.Length; $i++) { $ws.Cells.Item(2, $i+1).Value = $row2[$i] }

$row3 = @("ECs", "Efnb2", "Ephb4", "FAPs", 3, 1, 36.899643, 110.698929, 0.7238945645409351, 0.7238945645409351, 3, 1, 4.865208333333334, 14.595625, 0.1461844014571983, 0.1461844014571983, 179.524450620625, 1615.720055585625, 0.1058220936355358, 0.1058220936355358)
for ($i=0; $i -lt $row3.Length; $i++) { $ws.Cells.Item(3, $i+1).Value = $row3[$i] }

$row4 = @("ECs", "Efnb2", "Ephb4", "sCs", 3, 1, 36.899643, 110.698929, 0.7238945645409351, 0.7238945645409351, 3, 1, 4.669782333333333, 14.009347, 0.1403124570548501, 0.1403124570548501, 172.313300987707, 1550.819708889363, 0.1015714249993894, 0.1015714249993894)
for ($i=0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, $i+1).Value = $row4[$i] }

$row5 = @("FAPs", "Efnb2", "Ephb4", "ECs", 3, 1, 3.374819, 10.124457, 0.0662069584361419, 0.0662069584361419, 2, 0.6666666666666666, 23.746319, 71.238957, 0.7135031414879517, 0.7135031414879517, 80.139528541261, 721.2557568713489, 0.04723887283254949, 0.04723887283254949)
for ($i=0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, $i+1).Value = $row5[$i] }

$row6 = @("FAPs", "Efnb2", "Ephb4", "FAPs", 3, 1, 3.374819, 10.124457, 0.0662069584361419, 0.0662069584361419, 3, 1, 4.865208333333334, 14.595625, 0.1461844014571983, 0.1461844014571983, 16.41919752229167, 147.772777700625, 0.009678424591289009, 0.009678424591289007)
for ($i=0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, $i+1).Value = $row6[$i] }

$row7 = @("FAPs", "Efnb2", "Ephb4", "sCs", 3, 1, 3.374819, 10.124457, 0.0662069584361419, 0.0662069584361419, 3, 1, 4.669782333333333, 14.009347, 0.1403124570548501, 0.1403124570548501, 15.75967014439767, 141.837031299579, 0.009289661012303404, 0.009289661012303404)
for ($i=0; $i -lt $row7.Length; $i++) { $ws.Cells.Item(7, $i+1).Value = $row7[$i] }

$row8 = @("sCs", "Efnb2", "Ephb4", "ECs", 3, 1, 10.699319, 32.097957, 0.2098984770229228, 0.2098984770229228, 2, 0.6666666666666666, 23.746319, 71.238957, 0.7135031414879517, 0.7135031414879517, 254.069442056761, 2286.624978510849, 0.1497632227493921, 0.1497632227493921)
for ($i=0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, $i+1).Value = $row8[$i] }

$row9 = @("sCs", "Efnb2", "Ephb4", "FAPs", 3, 1, 10.699319, 32.097957, 0.2098984770229228, 0.2098984770229228, 3, 1, 4.865208333333334, 14.595625, 0.1461844014571983, 0.1461844014571983, 52.05441595979168, 468.489743638125, 0.03068388323037347, 0.03068388323037346)
for ($i=0; $i -lt $row9.Length; $i++) { $ws.Cells.Item(9, $i+1).Value = $row9[$i] }

$row10 = @("sCs", "Efnb2", "Ephb4", "sCs", 3, 1, 10.699319, 32.097957, 0.2098984770229228, 0.2098984770229228, 3, 1, 4.669782333333333, 14.009347, 0.1403124570548501, 0.1403124570548501, 49.96349084489767, 449.671417604079, 0.02945137104315729, 0.02945137104315729)
for ($i=0; $i -lt $row10.Length; $i++) { $ws.Cells.Item(10, $i+1).Value = $row10[$i] }

